$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(3.182878228561681, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 4.733082622659194)
    3 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    4 = @(0.000009318123435519965, 0.004309184025731883, 157.8057217802531, 6.48142807727062, 164.2914683596729)
    5 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 9.576116808119359)
    6 = @(0.1554434735375247, 0.3375848360084654, 157.8057217802531, 6.48142807727062, 164.7801781670697)
    7 = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 6.48142807727062, 26.62400969366105)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value2 = $vals[0]
    $ws.Range("C$row").Value2 = $vals[1]
    $ws.Range("D$row").Value2 = $vals[2]
    $ws.Range("E$row").Value2 = $vals[3]
    $ws.Range("G$row").Value2 = $vals[4]
}
